# Update "想去人数" (want-to-go count) figures for the two sheets that
# contain the conference data: "展览" (sheet 1) and "全部类型" (sheet 4).
# Rows 2, 3, 5 and 8 in column F each get bumped up by 2.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 71
    $ws.Range("F3").Value = 1432
    $ws.Range("F5").Value = 13
    $ws.Range("F8").Value = 209
}
